# Results: Update income-related prompts and calculations for percentage of
# housing and transport costs (fix: #5, fix: #18)
#
# - Clarify the household-income prompts (FR/EN) in the Widgets sheet.
# - Widen the two top household-income brackets (150k-199999/200k+ ->
#   150k-209999/210k+) in the Choices sheet.
# - Wrap the "percentage of gross income" label in parentheses in the
#   Labels sheet.
# - Re-point the workbook's active tab/selection at the Choices sheet, and
#   refresh the remembered selections on the Widgets/Choices/Labels sheets.

$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------------
# 1) Widgets sheet: household-income question prompts (FR/EN) + row height
# ---------------------------------------------------------------------
$wsWidgets = $wb.Worksheets.Item("Widgets")
$wsWidgets.Activate()

$wsWidgets.Range("G4").Value = "**Tranche de revenu** avant impôts (brut) **du ménage**, en 2024?`n__Cette information sert à calculer le pourcentage du revenu brut consacré aux transports et au logement.__"
$wsWidgets.Range("H4").Value = "What was your **household's income range** before taxes (gross income), in 2024?`n__This information is used to calculate the percentage of gross income spent on transportation and housing.__"

$wsWidgets.Rows.Item(4).RowHeight = 124.6

# Restore the frozen header pane (row 1 / column A) and move the
# remembered scroll / selection up to B5 / H4.
$wsWidgets.Range("B5").Select()
$excel.ActiveWindow.FreezePanes = $true
$wsWidgets.Range("H4").Select()

# ---------------------------------------------------------------------
# 2) Choices sheet: widen the top two household-income brackets
# ---------------------------------------------------------------------
$wsChoices = $wb.Worksheets.Item("Choices")
$wsChoices.Activate()

$wsChoices.Range("B33").Value = "150000_209999"
$wsChoices.Range("C33").Value = "150 000$ à 209 999$"
$wsChoices.Range("D33").Value = "$150,000 to $209,999"

$wsChoices.Range("B34").Value = "210000_999999"
$wsChoices.Range("C34").Value = "210 000$ et plus"
$wsChoices.Range("D34").Value = "$210,000 and more"

$wsChoices.Range("E32").Select()

# ---------------------------------------------------------------------
# 3) Labels sheet: parenthesize the percentage-of-income label (FR/EN)
# ---------------------------------------------------------------------
$wsLabels = $wb.Worksheets.Item("Labels")
$wsLabels.Activate()

$wsLabels.Range("C20").Value = "({{percentageOfIncome}}% du revenu brut)"
$wsLabels.Range("D20").Value = "({{percentageOfIncome}}% of gross income)"

$wsLabels.Range("C20").Select()

# ---------------------------------------------------------------------
# 4) Make Choices the workbook's active/selected tab
# ---------------------------------------------------------------------
$wsChoices.Activate()
